{"js": "// Apply the Arabic-translation edits described in the diff.\n// Each edit is done via a body.search() for the exact original text\n// (including any trailing space that belongs to the edited run) followed\n// by an insertText(..., Word.InsertLocation.replace) with the new text.\n\nasync function replaceOnce(context, searchText, newText, matchCase) {\n  const results = context.document.body.search(searchText, {\n    matchCase: matchCase !== false,\n    matchWildcards: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"{Onboarding} \" -> \"{Onboarding}\" (title heading, drop trailing space)\nawait replaceOnce(context, \"{Onboarding} \", \"{Onboarding}\");\n\n// 2. \"{Lesson: Onboarding} \" -> \"{Lesson: Onboarding}\" (drop trailing space)\nawait replaceOnce(context, \"{Lesson: Onboarding} \", \"{Lesson: Onboarding}\");\n\n// 3. \" \u0639\u0645\u0644. \" -> \" \u0639\u0645\u0644\u0647\u0627. \"\nawait replaceOnce(context, \" \u0639\u0645\u0644. \", \" \u0639\u0645\u0644\u0647\u0627. \");\n\n// 4. Replace the English \"Review Tips\" quoted phrase with its Arabic\n//    translation, and tidy up the extra/trailing spaces around it.\nawait replaceOnce(\n  context,\n  \"\u0625\u0630\u0627 \u0643\u0646\u062a \u062a\u0631\u063a\u0628 \u0641\u064a \u0645\u0631\u0627\u062c\u0639\u0629 \u0623\u064a \u0645\u0646 \u0627\u0644\u0646\u0635\u0627\u0626\u062d \u0627\u0644\u062a\u064a \u062a\u0644\u0642\u064a\u062a\u0647\u0627 \u0633\u0627\u0628\u0642\u064b\u0627\u060c \u0645\u0627 \u0639\u0644\u064a\u0643 \u0633\u0648\u0649 \u0643\u062a\u0627\u0628\u0629 MENU \u0648\u0627\u0644\u0627\u0646\u062a\u0642\u0627\u0644 \u0625\u0644\u0649   \\u201cReview Tips\\u201d \",\n  \"\u0625\u0630\u0627 \u0643\u0646\u062a \u062a\u0631\u063a\u0628 \u0641\u064a \u0645\u0631\u0627\u062c\u0639\u0629 \u0623\u064a \u0645\u0646 \u0627\u0644\u0646\u0635\u0627\u0626\u062d \u0627\u0644\u062a\u064a \u062a\u0644\u0642\u064a\u062a\u0647\u0627 \u0633\u0627\u0628\u0642\u064b\u0627\u060c \u0645\u0627 \u0639\u0644\u064a\u0643 \u0633\u0648\u0649 \u0643\u062a\u0627\u0628\u0629 MENU \u0648\u0627\u0644\u0627\u0646\u062a\u0642\u0627\u0644 \u0625\u0644\u0649 \\u201c\u0646\u0635\u0627\u0626\u062d \u0644\u0644\u0645\u0631\u0627\u062c\u0639\u0629\\u201d\"\n);\n\n// 5. Replace the whole (duplicate/old) paragraph with new copy about\n//    changing language/gender settings.\nawait replaceOnce(\n  context,\n  \"\u0625\u0630\u0627 \u0643\u0646\u062a \u062a\u0631\u063a\u0628 \u0641\u064a \u0645\u0631\u0627\u062c\u0639\u0629 \u0623\u064a \u0645\u0646 \u0627\u0644\u0646\u0635\u0627\u0626\u062d \u0627\u0644\u062a\u064a \u062a\u0644\u0642\u064a\u062a\u0647\u0627 \u0633\u0627\u0628\u0642\u064b\u0627\u060c \u0645\u0627 \u0639\u0644\u064a\u0643 \u0633\u0648\u0649 \u0643\u062a\u0627\u0628\u0629 \\\"\u0627\u0644\u0642\u0627\u0626\u0645\u0629\\\"\u0648\u0627\u0644\u0627\u0646\u062a\u0642\u0627\u0644 \u0625\u0644\u0649 \\u201c\u0645\u0631\u0627\u062c\u0639\u0629 \u0627\u0644\u0646\u0635\u0627\u0626\u062d\\u201d\",\n  \"\u0644\u062a\u063a\u064a\u064a\u0631 \u0625\u0639\u062f\u0627\u062f\u0627\u062a \u0627\u0644\u0644\u063a\u0629 \u0623\u0648 \u0627\u0644\u062c\u0646\u0633\u060c \u0627\u062e\u062a\u0631 \\\"\u062a\u063a\u064a\u064a\u0631 \u0625\u0639\u062f\u0627\u062f\u0627\u062a\u064a\\\"\"\n);\n\n// 6. \"\u062f\u0639\u0648\u0629 \u0635\u062f\u064a\u0642 \u0644 \" -> \"\u0627\u062f\u0639\u064f \u0635\u062f\u064a\u0642\u064b\u0627 \u0644 \" (keep trailing space; field code follows)\nawait replaceOnce(\n  context,\n  \"\u0644\u0645\u0634\u0627\u0631\u0643\u0629 \u0631\u0627\u0628\u0637 \u0647\u0630\u0627 \u0627\u0644\u0631\u0648\u0628\u0648\u062a \u0645\u0639 \u0635\u062f\u064a\u0642\u060c \u0627\u062e\u062a\u0631 \\u201c\u062f\u0639\u0648\u0629 \u0635\u062f\u064a\u0642 \u0644 \",\n  \"\u0644\u0645\u0634\u0627\u0631\u0643\u0629 \u0631\u0627\u0628\u0637 \u0647\u0630\u0627 \u0627\u0644\u0631\u0648\u0628\u0648\u062a \u0645\u0639 \u0635\u062f\u064a\u0642\u060c \u0627\u062e\u062a\u0631 \\u201c\u0627\u062f\u0639\u064f \u0635\u062f\u064a\u0642\u064b\u0627 \u0644 \"\n);\n\n// 7. \"\u0623\u0647\u0644\u0627 \u0628\u0643 \u0628 \" -> \"\u0623\u0647\u0644\u0627 \u0628\u0643 \u0641\u064a\" (drop trailing space; field code follows)\nawait replaceOnce(\n  context,\n  \"\u0648\u062c\u0648\u062f\u0643 \u0647\u0646\u0627 \u064a\u0638\u0647\u0631 \u0623\u0646\u0643 \u062a\u0647\u062a\u0645. \u0623\u0647\u0644\u0627 \u0628\u0643 \u0628 \",\n  \"\u0648\u062c\u0648\u062f\u0643 \u0647\u0646\u0627 \u064a\u0638\u0647\u0631 \u0623\u0646\u0643 \u062a\u0647\u062a\u0645. \u0623\u0647\u0644\u0627 \u0628\u0643 \u0641\u064a\"\n);\n\n// 8. \" .\" -> \".\" (trailing period run after the \"{Programme Name}\" field\n//    in the welcome paragraph; drop the leading space)\nawait replaceOnce(context, \" .\", \".\");\n", "ps1": "# Apply the Arabic-translation edits described in the diff.\n#\n# Each edit: locate the exact original run text with Find.Execute (no\n# Replace argument, so Word's \"smart quotes\" AutoFormat never touches the\n# text we supply) and then set .Text directly on the found Range \u2014 this\n# keeps straight quotes as straight quotes for the paragraph that needs\n# them (change 5 below).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n    if (-not $ok) {\n        throw \"Text not found: $findText\"\n    }\n    $rng.Text = $replaceText\n}\n\n# 1. \"{Onboarding} \" -> \"{Onboarding}\" (title heading, drop trailing space)\nReplace-Text \"{Onboarding} \" \"{Onboarding}\"\n\n# 2. \"{Lesson: Onboarding} \" -> \"{Lesson: Onboarding}\" (drop trailing space)\nReplace-Text \"{Lesson: Onboarding} \" \"{Lesson: Onboarding}\"\n\n# 3. \" \u0639\u0645\u0644. \" -> \" \u0639\u0645\u0644\u0647\u0627. \"\nReplace-Text \" \u0639\u0645\u0644. \" \" \u0639\u0645\u0644\u0647\u0627. \"\n\n# 4. Replace the English \"Review Tips\" quoted phrase with its Arabic\n#    translation, and tidy up the extra/trailing spaces around it.\nReplace-Text \"\u0625\u0630\u0627 \u0643\u0646\u062a \u062a\u0631\u063a\u0628 \u0641\u064a \u0645\u0631\u0627\u062c\u0639\u0629 \u0623\u064a \u0645\u0646 \u0627\u0644\u0646\u0635\u0627\u0626\u062d \u0627\u0644\u062a\u064a \u062a\u0644\u0642\u064a\u062a\u0647\u0627 \u0633\u0627\u0628\u0642\u064b\u0627\u060c \u0645\u0627 \u0639\u0644\u064a\u0643 \u0633\u0648\u0649 \u0643\u062a\u0627\u0628\u0629 MENU \u0648\u0627\u0644\u0627\u0646\u062a\u0642\u0627\u0644 \u0625\u0644\u0649   \u201cReview Tips\u201d \" \"\u0625\u0630\u0627 \u0643\u0646\u062a \u062a\u0631\u063a\u0628 \u0641\u064a \u0645\u0631\u0627\u062c\u0639\u0629 \u0623\u064a \u0645\u0646 \u0627\u0644\u0646\u0635\u0627\u0626\u062d \u0627\u0644\u062a\u064a \u062a\u0644\u0642\u064a\u062a\u0647\u0627 \u0633\u0627\u0628\u0642\u064b\u0627\u060c \u0645\u0627 \u0639\u0644\u064a\u0643 \u0633\u0648\u0649 \u0643\u062a\u0627\u0628\u0629 MENU \u0648\u0627\u0644\u0627\u0646\u062a\u0642\u0627\u0644 \u0625\u0644\u0649 \u201c\u0646\u0635\u0627\u0626\u062d \u0644\u0644\u0645\u0631\u0627\u062c\u0639\u0629\u201d\"\n\n# 5. Replace the whole (duplicate/old) paragraph with new copy about\n#    changing language/gender settings.\nReplace-Text \"\u0625\u0630\u0627 \u0643\u0646\u062a \u062a\u0631\u063a\u0628 \u0641\u064a \u0645\u0631\u0627\u062c\u0639\u0629 \u0623\u064a \u0645\u0646 \u0627\u0644\u0646\u0635\u0627\u0626\u062d \u0627\u0644\u062a\u064a \u062a\u0644\u0642\u064a\u062a\u0647\u0627 \u0633\u0627\u0628\u0642\u064b\u0627\u060c \u0645\u0627 \u0639\u0644\u064a\u0643 \u0633\u0648\u0649 \u0643\u062a\u0627\u0628\u0629 \"\"\u0627\u0644\u0642\u0627\u0626\u0645\u0629\"\"\u0648\u0627\u0644\u0627\u0646\u062a\u0642\u0627\u0644 \u0625\u0644\u0649 \u201c\u0645\u0631\u0627\u062c\u0639\u0629 \u0627\u0644\u0646\u0635\u0627\u0626\u062d\u201d\" \"\u0644\u062a\u063a\u064a\u064a\u0631 \u0625\u0639\u062f\u0627\u062f\u0627\u062a \u0627\u0644\u0644\u063a\u0629 \u0623\u0648 \u0627\u0644\u062c\u0646\u0633\u060c \u0627\u062e\u062a\u0631 \"\"\u062a\u063a\u064a\u064a\u0631 \u0625\u0639\u062f\u0627\u062f\u0627\u062a\u064a\"\"\"\n\n# 6. \"\u062f\u0639\u0648\u0629 \u0635\u062f\u064a\u0642 \u0644 \" -> \"\u0627\u062f\u0639\u064f \u0635\u062f\u064a\u0642\u064b\u0627 \u0644 \" (keep trailing space; field code follows)\nReplace-Text \"\u0644\u0645\u0634\u0627\u0631\u0643\u0629 \u0631\u0627\u0628\u0637 \u0647\u0630\u0627 \u0627\u0644\u0631\u0648\u0628\u0648\u062a \u0645\u0639 \u0635\u062f\u064a\u0642\u060c \u0627\u062e\u062a\u0631 \u201c\u062f\u0639\u0648\u0629 \u0635\u062f\u064a\u0642 \u0644 \" \"\u0644\u0645\u0634\u0627\u0631\u0643\u0629 \u0631\u0627\u0628\u0637 \u0647\u0630\u0627 \u0627\u0644\u0631\u0648\u0628\u0648\u062a \u0645\u0639 \u0635\u062f\u064a\u0642\u060c \u0627\u062e\u062a\u0631 \u201c\u0627\u062f\u0639\u064f \u0635\u062f\u064a\u0642\u064b\u0627 \u0644 \"\n\n# 7. \"\u0623\u0647\u0644\u0627 \u0628\u0643 \u0628 \" -> \"\u0623\u0647\u0644\u0627 \u0628\u0643 \u0641\u064a\" (drop trailing space; field code follows)\nReplace-Text \"\u0648\u062c\u0648\u062f\u0643 \u0647\u0646\u0627 \u064a\u0638\u0647\u0631 \u0623\u0646\u0643 \u062a\u0647\u062a\u0645. \u0623\u0647\u0644\u0627 \u0628\u0643 \u0628 \" \"\u0648\u062c\u0648\u062f\u0643 \u0647\u0646\u0627 \u064a\u0638\u0647\u0631 \u0623\u0646\u0643 \u062a\u0647\u062a\u0645. \u0623\u0647\u0644\u0627 \u0628\u0643 \u0641\u064a\"\n\n# 8. \" .\" -> \".\" (trailing period run after the \"{Programme Name}\" field\n#    in the welcome paragraph; drop the leading space)\nReplace-Text \" .\" \".\"\n"}
